# Updated results for 12/11/24
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New game result row (left table, columns A:F) ---
# Row 30: Hawks @ Knicks, 12/11/24
$ws.Range("A30").Value = 45637
$ws.Range("B30").Value = "Hawks"
$ws.Range("C30").Value = "Knicks"
$ws.Range("D30").Value = 237.5
$ws.Range("E30").Value = "L"
$ws.Range("F30").Formula = "=IF(E30=""L"",-1,IF(E30=""W"",1/1.1,0))"
$ws.Range("A29:F29").Copy()
$ws.Range("A30:F30").PasteSpecial(-4122)

# --- New bet row (right table, columns J:P) ---
# Row 34: Warriors @ Rockets, 12/11/24, "Warriors +2.5"
$ws.Range("J34").Value = 45637
$ws.Range("K34").Value = "Warriors"
$ws.Range("L34").Value = "Rockets"
$ws.Range("M34").Value = "Warriors +2.5"
$ws.Range("N34").Value = -110
$ws.Range("O34").Value = "W"
$ws.Range("P34").Formula = "=IF(O34=""L"",-1,IF(N34<0,1/(-N34/100),1*(N34/100)))"
$ws.Range("J33:P33").Copy()
$ws.Range("J34:P34").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update selection to match the saved workbook state
[void]$ws.Range("J35").Select()
